$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that must stay as Text (the source data uses
# "." as a thousands separator and keeps trailing/leading zeros), so force the
# Text number format before writing - otherwise Excel would auto-convert
# numeric-looking values (e.g. "578.10", "0.0000180") into actual numbers and
# silently reformat them (dropping the trailing zero, switching to scientific
# notation, etc).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '63.151.98'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '2.477.07'
$ws.Range("E3").Value = '  +2.46%  '
$ws.Range("E4").Value = '  -0.83%  '
$ws.Range("D5").Value = '578.10'
$ws.Range("E5").Value = '  +0.89%  '
$ws.Range("D6").Value = '147.10'
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  -0.36%  '
$ws.Range("D9").Value = '2.475.77'
$ws.Range("E9").Value = '  +1.12%  '
$ws.Range("E10").Value = '  +0.51%  '
$ws.Range("E11").Value = '  +1.60%  '
$ws.Range("E12").Value = '  +0.89%  '
$ws.Range("D13").Value = '0.354'
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("D14").Value = '28.72'
$ws.Range("E14").Value = '  +4.81%  '
$ws.Range("D15").Value = '0.0000180'
$ws.Range("E15").Value = '  +0.98%  '
$ws.Range("D16").Value = '2.926.43'
$ws.Range("E16").Value = '  +2.39%  '
$ws.Range("D17").Value = '63.091.83'
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("D18").Value = '2.473.58'
$ws.Range("E18").Value = '  +0.97%  '
$ws.Range("D19").Value = '8.23'
$ws.Range("E19").Value = '  +4.44%  '
$ws.Range("D20").Value = '11.08'
$ws.Range("E20").Value = '  +0.77%  '
$ws.Range("D21").Value = '329.55'
$ws.Range("E21").Value = '  +0.27%  '
$ws.Range("D22").Value = '2.26'
$ws.Range("E22").Value = '  +9.30%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").Value = '66.39'
$ws.Range("E25").Value = '  +1.27%  '
$ws.Range("D26").Value = '672.52'
$ws.Range("E26").Value = '  +5.84%  '
$ws.Range("D27").Value = '9.69'
$ws.Range("E27").Value = '  +13.05%  '
$ws.Range("D28").Value = '0.0₃0998'
$ws.Range("E28").Value = '  +0.96%  '
$ws.Range("D29").Value = '2.594.53'
$ws.Range("E29").Value = '  +2.14%  '
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  +402.30%  '
$ws.Range("E31").Value = '  +2.76%  '
$ws.Range("D32").Value = '8.08'
$ws.Range("E32").Value = '  -1.91%  '
$ws.Range("E33").Value = '  +1.61%  '
$ws.Range("E34").Value = '  -4.04%  '
$ws.Range("D35").Value = '1.55'
$ws.Range("E35").Value = '  +3.87%  '
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  +0.40%  '
$ws.Range("E37").Value = '  +0.75%  '
$ws.Range("E38").Value = '  +0.95%  '
$ws.Range("E39").Value = '  -0.56%  '
$ws.Range("D40").Value = '18.82'
$ws.Range("E40").Value = '  +0.72%  '
$ws.Range("D41").Value = '151.97'
$ws.Range("E41").Value = '  -0.71%  '
$ws.Range("D42").Value = '2.74'
$ws.Range("E42").Value = '  -0.88%  '
$ws.Range("E43").Value = '  -0.42%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = '0.0₆0310'
$ws.Range("E45").Value = '  +8.50%  '
$ws.Range("D46").Value = '154.28'
$ws.Range("E46").Value = '  +6.21%  '
$ws.Range("E47").Value = '  +19.35%  '
$ws.Range("D48").Value = '3.62'
$ws.Range("E48").Value = '  +0.38%  '
$ws.Range("D49").Value = '20.69'
$ws.Range("E49").Value = '  +0.92%  '
$ws.Range("D50").Value = '0.607'
$ws.Range("E50").Value = '  +0.81%  '
$ws.Range("E51").Value = '  -0.89%  '
